$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Figure out how many data rows currently exist (column A holds the
# segment name in each row, starting at row 2).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
if ($lastRow -lt 2) { $lastRow = 1 }

# Stash the existing segment names (currently in column A) before the
# column insert shuffles everything one slot to the right.
$segmentNames = @()
for ($row = 2; $row -le $lastRow; $row++) {
    $segmentNames += $ws.Cells.Item($row, 1).Value()
}

# Insert a new column before column B; existing B:K shift right to C:L.
$ws.Columns.Item(2).Insert()

# Give the new header cell (B1) the same style as the other header cells.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("B1").Value = "segments"

# Move the segment names into the new column B, and replace column A
# with a plain 0-based numeric index.
for ($i = 0; $i -lt $segmentNames.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i
    $ws.Cells.Item($row, 2).Value = $segmentNames[$i]
}

# The new segment-name column keeps no special formatting (unlike the old
# bordered/bold/centered column A), so clear what Insert copied over.
$ws.Range("B2:B" + $lastRow).ClearFormats()

Write-Output "done"
